# Add a new "2022-Q3" sheet (cloned from "2022-Q2"), insert it right after
# "总计", refresh its holdings figures, and record the new quarter in the
# "总计" summary sheet (shifting the older rows down by one).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Create the "2022-Q3" sheet by copying "2022-Q2" and moving the copy
#    right after "总计" (tab position 2), then rename it.
# ---------------------------------------------------------------------
$sourceSheet = $wb.Worksheets.Item("2022-Q2")
$afterSheet = $wb.Worksheets.Item("总计")
$sourceSheet.Copy($null, $afterSheet)

$q3 = $wb.Worksheets.Item(2)
$q3.Name = "2022-Q3"

# ---------------------------------------------------------------------
# 2) Refresh the 2022-Q3 holdings figures (fund codes/names stay the same
#    as the prior quarter; size/position/value/rank are updated). The
#    D:G columns hold text-formatted numbers (e.g. "15.90"), so force the
#    Text number format first to keep them from being coerced to numeric.
# ---------------------------------------------------------------------
$textCols = $q3.Range("D2:G4")
$textCols.NumberFormat = "@"

$q3.Range("D2").Value = "15.90"
$q3.Range("E2").Value = "99.95"
$q3.Range("F2").Value = "4.21"
$q3.Range("G2").Value = "0.6694"
$q3.Range("H2").Value = 6

$q3.Range("D3").Value = "14.02"
$q3.Range("E3").Value = "94.15"
$q3.Range("F3").Value = "4.27"
$q3.Range("G3").Value = "0.5987"
$q3.Range("H3").Value = 7

$q3.Range("D4").Value = "0.39"
$q3.Range("E4").Value = "94.15"
$q3.Range("F4").Value = "4.27"
$q3.Range("G4").Value = "0.0167"
$q3.Range("H4").Value = 7

# ---------------------------------------------------------------------
# 3) Update the "总计" summary sheet: push existing data rows down by one
#    and insert the new 2022-Q3 totals at the top.
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")
$total.Rows.Item(2).Insert()

$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q3"
$total.Range("C2").Value = 3
$total.Range("D2").Value = 1.28

# Renumber the index column (A) for the rows that shifted down.
$total.Range("A3").Value = 1
$total.Range("A4").Value = 2
$total.Range("A5").Value = 3
$total.Range("A6").Value = 4
$total.Range("A7").Value = 5
$total.Range("A8").Value = 6
